$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column so numeric-looking strings
# (e.g. "1.00", "61.08") are preserved exactly as text rather than being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.961.32"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.342.47"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "517.87"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "135.57"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "2.354.44"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +5.26%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "23.97"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "2.761.31"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "56.950.81"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "2.356.57"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "327.17"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "6.76"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "61.08"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "0.165"
$ws.Range("E25").Value = "  +4.84%  "
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "7.98"
$ws.Range("E27").Value = "  +4.04%  "
$ws.Range("D28").Value = "1.29"
$ws.Range("E28").Value = "  +9.32%  "
$ws.Range("D29").Value = "170.43"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "18.57"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "0.917"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "4.02"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("D41").Value = "147.48"
$ws.Range("E41").Value = "  +7.05%  "
$ws.Range("D42").Value = "0.384"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "3.64"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "281.12"
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("D45").Value = "5.24"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "0.0936"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "0.0507"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "0.563"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "18.01"
$ws.Range("E50").Value = "  +6.24%  "
$ws.Range("B51").Value = "Polygon"
$ws.Range("C51").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D51").Value = "0.383"
$ws.Range("E51").Value = "  -0.04%  "

# Restore the default cell style so the text-format override above
# does not leave a lingering explicit number format on the cells.
$ws.Range("D2:D51").Style = "Normal"

